$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 4143565.5
$ws.Range("I18").Value = 991.8
$ws.Range("K18").Value = 991.8
$ws.Range("M18").Value = -707.8

$ws.Range("H32").Value = 8999
$ws.Range("J32").Value = 8999
$ws.Range("L32").Value = 8999
$ws.Range("N32").Value = -9651

$ws.Range("H33").Value = 850.72
$ws.Range("I33").Value = 280.57144
$ws.Range("K33").Value = 280.57144
$ws.Range("M33").Value = -51.57144

$ws.Range("H34").Value = 23925
$ws.Range("I34").Value = 23925
$ws.Range("K34").Value = 23925
$ws.Range("M34").Value = -23722

$ws.Range("H36").Value = 23925
$ws.Range("I36").Value = 23925
$ws.Range("K36").Value = 23925
$ws.Range("M36").Value = -23210

$ws.Range("H40").Value = 4054.625
$ws.Range("J40").Value = 6000
$ws.Range("L40").Value = 6000
$ws.Range("N40").Value = -6350

$ws.Range("H58").Value = 1524
$ws.Range("J58").Value = 2500
$ws.Range("L58").Value = 7500
$ws.Range("N58").Value = -7800

$ws.Range("H69").Value = 16568.725
$ws.Range("I69").Value = 11111.777
$ws.Range("J69").Value = 19024.35
$ws.Range("K69").Value = 33335.331
$ws.Range("L69").Value = 57073.05
$ws.Range("M69").Value = -32461.331
$ws.Range("N69").Value = -58821.05

$ws.Range("H70").Value = 3521.7144
$ws.Range("I70").Value = 2160.4
$ws.Range("J70").Value = 6925
$ws.Range("K70").Value = 6481.200000000001
$ws.Range("L70").Value = 20775
$ws.Range("M70").Value = -6211.200000000001
$ws.Range("N70").Value = -21315

$ws.Range("H72").Value = 16568.725
$ws.Range("I72").Value = 11111.777
$ws.Range("J72").Value = 19024.35
$ws.Range("K72").Value = 100005.993
$ws.Range("L72").Value = 171219.15
$ws.Range("M72").Value = -95637.993
$ws.Range("N72").Value = -179955.15

$ws.Range("H73").Value = 3521.7144
$ws.Range("I73").Value = 2160.4
$ws.Range("J73").Value = 6925
$ws.Range("K73").Value = 6481.200000000001
$ws.Range("L73").Value = 20775
$ws.Range("M73").Value = -5545.200000000001
$ws.Range("N73").Value = -22647

$ws.Range("H92").Value = 1725
$ws.Range("J92").Value = 1699.5
$ws.Range("L92").Value = 1699.5
$ws.Range("N92").Value = -4195.5

$ws.Range("H98").Value = 1668148.6
$ws.Range("I98").Value = 1895410
$ws.Range("K98").Value = 1895410
$ws.Range("M98").Value = -1893912

$ws.Range("H100").Value = 911.44446
$ws.Range("I100").Value = 932.875
$ws.Range("J100").Value = 740
$ws.Range("K100").Value = 932.875
$ws.Range("L100").Value = 740
$ws.Range("M100").Value = -391.875
$ws.Range("N100").Value = -1822

$ws.Range("H112").Value = 3221.8157
$ws.Range("I112").Value = 1446.8889
$ws.Range("J112").Value = 3772.6553
$ws.Range("K112").Value = 4340.6667
$ws.Range("L112").Value = 11317.9659
$ws.Range("M112").Value = -3232.6667
$ws.Range("N112").Value = -13533.9659

$ws.Range("H122").Value = 1668148.6
$ws.Range("I122").Value = 1895410
$ws.Range("K122").Value = 5686230
$ws.Range("M122").Value = -5683780

$ws.Range("H132").Value = 2041.2368
$ws.Range("I132").Value = 1916.8125
$ws.Range("J132").Value = 2704.8333
$ws.Range("K132").Value = 5750.4375
$ws.Range("L132").Value = 8114.499899999999
$ws.Range("M132").Value = -3220.4375
$ws.Range("N132").Value = -13174.4999

$ws.Range("H135").Value = 1172
$ws.Range("J135").Value = 2999.6667
$ws.Range("L135").Value = 26997.0003
$ws.Range("N135").Value = -32067.0003

$ws.Range("H137").Value = 2488.879
$ws.Range("I137").Value = 2118.125
$ws.Range("K137").Value = 6354.375
$ws.Range("M137").Value = -3804.375

$ws.Range("H138").Value = 2322.88
$ws.Range("I138").Value = 974.4878
$ws.Range("J138").Value = 3259.8982
$ws.Range("K138").Value = 2923.4634
$ws.Range("L138").Value = 9779.694600000001
$ws.Range("M138").Value = 2216.5366
$ws.Range("N138").Value = -20059.6946

$ws.Range("H141").Value = 2136.0208
$ws.Range("I141").Value = 1036.4048
$ws.Range("J141").Value = 9833.333000000001
$ws.Range("K141").Value = 3109.2144
$ws.Range("L141").Value = 29499.999
$ws.Range("M141").Value = 2070.7856
$ws.Range("N141").Value = -39859.999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9213.530000000001
$ws.Range("I32").Value = 4477.5366
$ws.Range("K32").Value = 4477.5366
$ws.Range("M32").Value = -4190.5366

$ws.Range("H61").Value = 3250.1516
$ws.Range("I61").Value = 2395.3225
$ws.Range("J61").Value = 16500
$ws.Range("K61").Value = 2395.3225
$ws.Range("L61").Value = 16500
$ws.Range("M61").Value = -2183.3225
$ws.Range("N61").Value = -16924

$ws.Range("H63").Value = 3701.6667
$ws.Range("I63").Value = 2730.125
$ws.Range("K63").Value = 2730.125
$ws.Range("M63").Value = -2044.125

$ws.Range("H66").Value = 3701.6667
$ws.Range("I66").Value = 2730.125
$ws.Range("K66").Value = 13650.625
$ws.Range("M66").Value = -10218.625

$ws.Range("H74").Value = 4439.8213
$ws.Range("I74").Value = 2444.6843
$ws.Range("K74").Value = 2444.6843
$ws.Range("M74").Value = -1570.6843

$ws.Range("H77").Value = 4439.8213
$ws.Range("I77").Value = 2444.6843
$ws.Range("K77").Value = 12223.4215
$ws.Range("M77").Value = -7855.4215

$ws.Range("H97").Value = 4339.875
$ws.Range("I97").Value = 4817.143
$ws.Range("J97").Value = 999
$ws.Range("K97").Value = 4817.143
$ws.Range("L97").Value = 999
$ws.Range("M97").Value = -4321.143
$ws.Range("N97").Value = -1991

$ws.Range("H102").Value = 2629.625
$ws.Range("I102").Value = 2576.9285
$ws.Range("J102").Value = 2998.5
$ws.Range("K102").Value = 2576.9285
$ws.Range("L102").Value = 2998.5
$ws.Range("M102").Value = -954.9285
$ws.Range("N102").Value = -6242.5

$ws.Range("H122").Value = 2733.0732
$ws.Range("I122").Value = 2628.0303
$ws.Range("J122").Value = 3166.375
$ws.Range("K122").Value = 7884.090899999999
$ws.Range("L122").Value = 9499.125
$ws.Range("M122").Value = -5434.090899999999
$ws.Range("N122").Value = -14399.125

$ws.Range("H132").Value = 3907.6487
$ws.Range("I132").Value = 2612
$ws.Range("K132").Value = 7836
$ws.Range("M132").Value = -5306

$ws.Range("H136").Value = 3250.1516
$ws.Range("I136").Value = 2395.3225
$ws.Range("J136").Value = 16500
$ws.Range("K136").Value = 7185.967500000001
$ws.Range("L136").Value = 49500
$ws.Range("M136").Value = -4635.967500000001
$ws.Range("N136").Value = -54600

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H61").Value = 53999.5
$ws.Range("J61").Value = 53999.5
$ws.Range("L61").Value = 53999.5
$ws.Range("N61").Value = -54625.5

$ws.Range("H86").Value = 7578.6665
$ws.Range("I86").Value = 4231.4165
$ws.Range("J86").Value = 12041.667
$ws.Range("K86").Value = 4231.4165
$ws.Range("L86").Value = 12041.667
$ws.Range("M86").Value = -3108.4165
$ws.Range("N86").Value = -14287.667

$ws.Range("H89").Value = 7578.6665
$ws.Range("I89").Value = 4231.4165
$ws.Range("J89").Value = 12041.667
$ws.Range("K89").Value = 21157.0825
$ws.Range("L89").Value = 60208.335
$ws.Range("M89").Value = -15541.0825
$ws.Range("N89").Value = -71440.33499999999

$ws.Range("H94").Value = 767.1429000000001
$ws.Range("J94").Value = 2874
$ws.Range("L94").Value = 2874
$ws.Range("N94").Value = -3776

$ws.Range("H134").Value = 3996.35
$ws.Range("I134").Value = 3162.6667
$ws.Range("J134").Value = 11499.5
$ws.Range("K134").Value = 9488.000100000001
$ws.Range("L134").Value = 34498.5
$ws.Range("M134").Value = -6953.000100000001
$ws.Range("N134").Value = -39568.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 549.7143
$ws.Range("I7").Value = 415.57144
$ws.Range("J7").Value = 683.8570999999999
$ws.Range("K7").Value = 415.57144
$ws.Range("L7").Value = 683.8570999999999
$ws.Range("M7").Value = -302.57144
$ws.Range("N7").Value = -909.8570999999999

$ws.Range("H16").Value = 1136
$ws.Range("I16").Value = 1143.6666
$ws.Range("J16").Value = 1101.5
$ws.Range("K16").Value = 1143.6666
$ws.Range("L16").Value = 1101.5
$ws.Range("M16").Value = -856.6666
$ws.Range("N16").Value = -1675.5

$ws.Range("H22").Value = 1677.8667
$ws.Range("I22").Value = 2413.6667
$ws.Range("J22").Value = 574.1667
$ws.Range("K22").Value = 2413.6667
$ws.Range("L22").Value = 574.1667
$ws.Range("M22").Value = -2063.6667
$ws.Range("N22").Value = -1274.1667

$ws.Range("H31").Value = 5318.4375
$ws.Range("I31").Value = 2128.0344
$ws.Range("J31").Value = 10188
$ws.Range("K31").Value = 2128.0344
$ws.Range("L31").Value = 10188
$ws.Range("M31").Value = -1833.0344
$ws.Range("N31").Value = -10778

$ws.Range("H34").Value = 5318.4375
$ws.Range("I34").Value = 2128.0344
$ws.Range("J34").Value = 10188
$ws.Range("K34").Value = 2128.0344
$ws.Range("L34").Value = 10188
$ws.Range("M34").Value = -1926.0344
$ws.Range("N34").Value = -10592

$ws.Range("H58").Value = 2802.3447
$ws.Range("I58").Value = 2382.0417
$ws.Range("J58").Value = 4819.8
$ws.Range("K58").Value = 2382.0417
$ws.Range("L58").Value = 4819.8
$ws.Range("M58").Value = -2179.0417
$ws.Range("N58").Value = -5225.8

$ws.Range("H94").Value = 1879.1538
$ws.Range("I94").Value = 1807.6
$ws.Range("K94").Value = 1807.6
$ws.Range("M94").Value = -1356.6

$ws.Range("H113").Value = 1136
$ws.Range("I113").Value = 1143.6666
$ws.Range("J113").Value = 1101.5
$ws.Range("K113").Value = 1143.6666
$ws.Range("L113").Value = 1101.5
$ws.Range("M113").Value = 1026.3334
$ws.Range("N113").Value = -5441.5

$ws.Range("H122").Value = 2415.5557
$ws.Range("I122").Value = 2415.5557
$ws.Range("K122").Value = 7246.6671
$ws.Range("M122").Value = -4796.6671

$ws.Range("H132").Value = 2537.52
$ws.Range("I132").Value = 2068.0732
$ws.Range("J132").Value = 4676.1113
$ws.Range("K132").Value = 6204.219599999999
$ws.Range("L132").Value = 14028.3339
$ws.Range("M132").Value = -3674.219599999999
$ws.Range("N132").Value = -19088.3339

$ws.Range("H134").Value = 2806.9565
$ws.Range("I134").Value = 2229.375
$ws.Range("J134").Value = 4127.143
$ws.Range("K134").Value = 6688.125
$ws.Range("L134").Value = 12381.429
$ws.Range("M134").Value = -4153.125
$ws.Range("N134").Value = -17451.429

$ws.Range("H136").Value = 2802.3447
$ws.Range("I136").Value = 2382.0417
$ws.Range("J136").Value = 4819.8
$ws.Range("K136").Value = 7146.125100000001
$ws.Range("L136").Value = 14459.4
$ws.Range("M136").Value = -4596.125100000001
$ws.Range("N136").Value = -19559.4

$ws.Range("H141").Value = 199177.22
$ws.Range("J141").Value = 241799.42
$ws.Range("L141").Value = 241799.42
$ws.Range("N141").Value = -252159.42

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2418.889
$ws.Range("I5").Value = 2896
$ws.Range("J5").Value = 2180.3333
$ws.Range("K5").Value = 8688
$ws.Range("L5").Value = 6540.999899999999
$ws.Range("M5").Value = -8576
$ws.Range("N5").Value = -6764.999899999999

$ws.Range("H7").Value = 426.30768
$ws.Range("I7").Value = 352.6154
$ws.Range("J7").Value = 500
$ws.Range("K7").Value = 1057.8462
$ws.Range("L7").Value = 1500
$ws.Range("M7").Value = -945.8462
$ws.Range("N7").Value = -1724

$ws.Range("H23").Value = 126.27273
$ws.Range("I23").Value = 100
$ws.Range("J23").Value = 148.16667
$ws.Range("K23").Value = 300
$ws.Range("L23").Value = 444.50001
$ws.Range("M23").Value = -65
$ws.Range("N23").Value = -914.50001

$ws.Range("H26").Value = 286.25
$ws.Range("I26").Value = 286.25
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 858.75
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -570.75
$ws.Range("N26").Value = ""

$ws.Range("H132").Value = 1664
$ws.Range("J132").Value = 1500
$ws.Range("L132").Value = 13500
$ws.Range("N132").Value = -18560

$ws.Range("H135").Value = 2418.889
$ws.Range("I135").Value = 2896
$ws.Range("J135").Value = 2180.3333
$ws.Range("K135").Value = 26064
$ws.Range("L135").Value = 19622.9997
$ws.Range("M135").Value = -23529
$ws.Range("N135").Value = -24692.9997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 15781.429
$ws.Range("I49").Value = 4500
$ws.Range("J49").Value = 17661.666
$ws.Range("K49").Value = 4500
$ws.Range("L49").Value = 17661.666
$ws.Range("M49").Value = -4316
$ws.Range("N49").Value = -18029.666

$ws.Range("H80").Value = 2115.9
$ws.Range("J80").Value = 1876.25
$ws.Range("L80").Value = 1876.25
$ws.Range("N80").Value = -3872.25

$ws.Range("H83").Value = 2115.9
$ws.Range("J83").Value = 1876.25
$ws.Range("L83").Value = 9381.25
$ws.Range("N83").Value = -19365.25

$ws.Range("H97").Value = 719
$ws.Range("I97").Value = 549.5
$ws.Range("J97").Value = 912.7143
$ws.Range("K97").Value = 549.5
$ws.Range("L97").Value = 912.7143
$ws.Range("M97").Value = -53.5
$ws.Range("N97").Value = -1904.7143

$ws.Range("H102").Value = 3254.8206
$ws.Range("I102").Value = 2186.76
$ws.Range("J102").Value = 5162.0713
$ws.Range("K102").Value = 2186.76
$ws.Range("L102").Value = 5162.0713
$ws.Range("M102").Value = -564.7600000000002
$ws.Range("N102").Value = -8406.0713

$ws.Range("H122").Value = 7595.2905
$ws.Range("I122").Value = 4860.8945
$ws.Range("J122").Value = 11924.75
$ws.Range("K122").Value = 14582.6835
$ws.Range("L122").Value = 35774.25
$ws.Range("M122").Value = -12132.6835
$ws.Range("N122").Value = -40674.25

$ws.Range("H132").Value = 4571.787
$ws.Range("I132").Value = 4734.087
$ws.Range("J132").Value = 4416.25
$ws.Range("K132").Value = 14202.261
$ws.Range("L132").Value = 13248.75
$ws.Range("M132").Value = -11672.261
$ws.Range("N132").Value = -18308.75

$ws.Range("H141").Value = 27000
$ws.Range("J141").Value = 27000
$ws.Range("L141").Value = 27000
$ws.Range("N141").Value = -37360

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2695.889
$ws.Range("J46").Value = 3270.7144
$ws.Range("L46").Value = 3270.7144
$ws.Range("N46").Value = -3646.7144

$ws.Range("H100").Value = 2886.2222
$ws.Range("I100").Value = 3165.5
$ws.Range("K100").Value = 3165.5
$ws.Range("M100").Value = -2624.5

$ws.Range("H132").Value = 2910.5925
$ws.Range("I132").Value = 1560.081
$ws.Range("J132").Value = 5849.9414
$ws.Range("K132").Value = 4680.242999999999
$ws.Range("L132").Value = 17549.8242
$ws.Range("M132").Value = -2150.242999999999
$ws.Range("N132").Value = -22609.8242

$ws.Range("H134").Value = 60000
$ws.Range("J134").Value = 60000
$ws.Range("L134").Value = 60000
$ws.Range("N134").Value = -70140

$ws.Range("H136").Value = 7672.5063
$ws.Range("I136").Value = 4139.4614
$ws.Range("J136").Value = 9342.673000000001
$ws.Range("K136").Value = 12418.3842
$ws.Range("L136").Value = 28028.019
$ws.Range("M136").Value = -9868.3842
$ws.Range("N136").Value = -33128.019

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 82749.75
$ws.Range("J42").Value = 82749.75
$ws.Range("L42").Value = 82749.75
$ws.Range("N42").Value = -83505.75

$ws.Range("H70").Value = 9999
$ws.Range("J70").Value = 9999
$ws.Range("L70").Value = 9999
$ws.Range("N70").Value = -10629

$ws.Range("H73").Value = 9999
$ws.Range("J73").Value = 9999
$ws.Range("L73").Value = 9999
$ws.Range("N73").Value = -12183

$ws.Range("H97").Value = 145000
$ws.Range("J97").Value = 145000
$ws.Range("L97").Value = 145000
$ws.Range("N97").Value = -146982

$ws.Range("H107").Value = 1035.4062
$ws.Range("I107").Value = 1001.2692
$ws.Range("J107").Value = 1183.3334
$ws.Range("K107").Value = 3003.8076
$ws.Range("L107").Value = 3550.0002
$ws.Range("M107").Value = -1083.8076
$ws.Range("N107").Value = -7390.0002

$ws.Range("H122").Value = 2584.5625
$ws.Range("I122").Value = 2605.3928
$ws.Range("J122").Value = 2438.75
$ws.Range("K122").Value = 7816.178400000001
$ws.Range("L122").Value = 7316.25
$ws.Range("M122").Value = -5366.178400000001
$ws.Range("N122").Value = -12216.25

$ws.Range("H132").Value = 2574.377
$ws.Range("I132").Value = 2263.3022
$ws.Range("J132").Value = 3317.5
$ws.Range("K132").Value = 6789.9066
$ws.Range("L132").Value = 9952.5
$ws.Range("M132").Value = -4259.9066
$ws.Range("N132").Value = -15012.5

$ws.Range("H136").Value = 5808.6665
$ws.Range("I136").Value = 4924.1
$ws.Range("K136").Value = 14772.3
$ws.Range("M136").Value = -12222.3

$ws.Range("H140").Value = 29666.666
$ws.Range("I140").Value = 19500
$ws.Range("J140").Value = 50000
$ws.Range("K140").Value = 19500
$ws.Range("L140").Value = 50000
$ws.Range("M140").Value = -14320
$ws.Range("N140").Value = -60360
